$wb = $excel.ActiveWorkbook

# Use an existing header cell as the style source for the new header row
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Copy the header style (bold font, border, centered) from "Player Info"!A1:D1 onto the new header row
$ws1.Range("A1:D1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Match the page margins used by the other sheets in the workbook (in points: 72 pts = 1 in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row values
$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "BATTING_POSITION"
$ws.Cells.Item(1, 3).Value = "NUM_4"
$ws.Cells.Item(1, 4).Value = "NUM_6"
$ws.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Data row (text values are entered with a leading apostrophe so numeric-looking
# strings like "4625", "3", "0" and "10.42%" are kept as text instead of being
# converted to numbers/percentages; BATTING_POSITION is a genuine number)
$ws.Cells.Item(2, 1).Value = "'4625"
$ws.Cells.Item(2, 2).Value = 3
$ws.Cells.Item(2, 3).Value = "'3"
$ws.Cells.Item(2, 4).Value = "'0"
$ws.Cells.Item(2, 5).Value = "'10.42%"
$ws.Cells.Item(2, 6).Value = "'NO"
